# "Præver at slippe arket" - rework Ark1: drop the old A1:C3 number grid
# and lay down the new Vejleder1 / Vejleder2 / Møder header row at F3:H3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 3x3 block of numbers that used to live at A1:C3.
$ws.Range("A1:C3").ClearContents()

# Write the new header labels into F3:G3:H3 (these become shared strings).
$ws.Range("F3").Value = "Vejleder1"
$ws.Range("G3").Value = "Vejleder2"
$ws.Range("H3").Value = "Møder"

# Match the author's final selection/active cell (J4).
$ws.Range("J4").Select()
